$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily "Dolar observado" rows appended to the series.
$dates  = @("10-09-2021", "13-09-2021", "14-09-2021", "15-09-2021", "16-09-2021", "20-09-2021")
$values = @(791.28, 789.91, 784.26, 783.25, 781.85, 780.59)

$startRow = 176
$scratchRow = 500

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Write the date text through a throw-away helper cell using a text
    # formula so Excel's "typed by the user" date auto-detection never
    # kicks in (it would otherwise turn e.g. 10-09-2021 into a serial
    # date). Copy/PasteSpecial values-only into the real destination so
    # the destination cell keeps its original (default) style, then
    # remove the helper row again.
    $scratch = $ws.Cells.Item($scratchRow, 1)
    $scratch.Formula = '="' + $dates[$i] + '"'
    $scratch.Copy()
    $dest = $ws.Cells.Item($row, 1)
    $dest.PasteSpecial(-4163)
    $scratch.EntireRow.Delete()

    $ws.Cells.Item($row, 2).Value = $values[$i]
}
